$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 301
$ws.Range("F5").Value = 3043
$ws.Range("F9").Value = 301
$ws.Range("F10").Value = 7164
$ws.Range("F11").Value = 55
$ws.Range("F14").Value = 643
$ws.Range("F15").Value = 1552
$ws.Range("F16").Value = 2313
$ws.Range("F17").Value = 1568
$ws.Range("F18").Value = 1165
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 224
$ws.Range("F21").Value = 375
$ws.Range("F24").Value = 1851
$ws.Range("F27").Value = 51
$ws.Range("F28").Value = 1708
$ws.Range("F29").Value = 1299
$ws.Range("F30").Value = 156
$ws.Range("F32").Value = 27
$ws.Range("F34").Value = 471
$ws.Range("F35").Value = 46
$ws.Range("F36").Value = 2569
$ws.Range("F37").Value = 2862
$ws.Range("F38").Value = 2110
$ws.Range("F39").Value = 71
$ws.Range("F40").Value = 203
$ws.Range("F45").Value = 351
$ws.Range("F46").Value = 136
$ws.Range("F47").Value = 203
$ws.Range("F49").Value = 74

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 196
$ws.Range("F11").Value = 77
$ws.Range("F16").Value = 82
$ws.Range("F18").Value = 345
$ws.Range("F19").Value = 505
$ws.Range("F25").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1766
$ws.Range("F8").Value = 2834
$ws.Range("F9").Value = 1080
$ws.Range("F10").Value = 1025
$ws.Range("F12").Value = 368
$ws.Range("F13").Value = 1720
$ws.Range("F14").Value = 7839

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 301
$ws.Range("F4").Value = 3043
$ws.Range("F6").Value = 1766
$ws.Range("F7").Value = 301
$ws.Range("F8").Value = 2834
$ws.Range("F9").Value = 7164
$ws.Range("F10").Value = 1080
$ws.Range("F11").Value = 55
$ws.Range("F13").Value = 368
$ws.Range("F14").Value = 643
$ws.Range("F15").Value = 1552
$ws.Range("F16").Value = 2313
$ws.Range("F17").Value = 1568
$ws.Range("F18").Value = 1165
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 375
$ws.Range("F23").Value = 1851
$ws.Range("F24").Value = 77
$ws.Range("F26").Value = 51
$ws.Range("F27").Value = 1708
$ws.Range("F28").Value = 1299
$ws.Range("F29").Value = 156
$ws.Range("F31").Value = 27
$ws.Range("F33").Value = 82
$ws.Range("F35").Value = 505
$ws.Range("F36").Value = 471
$ws.Range("F37").Value = 46
$ws.Range("F38").Value = 2569
$ws.Range("F39").Value = 2862
$ws.Range("F40").Value = 2110
$ws.Range("F41").Value = 71
$ws.Range("F42").Value = 203
$ws.Range("F46").Value = 21
$ws.Range("F47").Value = 203
